$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new monthly row (row 45) that the update adds below the
# existing data (row 44 was the last one, for 01-07-2021).
#
# Column A holds a date-like text label ("01-08-2021"). Assigning that
# string straight to .Value makes Excel auto-detect it as a real date,
# which reformats the cell with a date number format (and drags a new,
# unwanted style into styles.xml). To keep it a plain text/shared-string
# cell - like all the other cells in column A - we first write it as a
# text-literal formula, then copy/paste-special the cell as values so the
# formula collapses down to its literal cached string.
$ws.Range("A45").Formula = "=""01-08-2021"""
$ws.Range("A45").Copy()
$ws.Range("A45").PasteSpecial(-4163)

$ws.Cells.Item(45, 2).Value = 114.79
$ws.Cells.Item(45, 3).Value = 108.39
$ws.Cells.Item(45, 4).Value = 98.64
$ws.Cells.Item(45, 5).Value = 112.08
$ws.Cells.Item(45, 6).Value = 111.59
$ws.Cells.Item(45, 7).Value = 107.16
$ws.Cells.Item(45, 8).Value = 111.31
$ws.Cells.Item(45, 9).Value = 93.40000000000001
$ws.Cells.Item(45, 10).Value = 109.06
$ws.Cells.Item(45, 11).Value = 112.29
$ws.Cells.Item(45, 12).Value = 110.96
$ws.Cells.Item(45, 13).Value = 111.8
